# Generate Report for Handback
# Updates the localization-status report: the two files that were
# previously "Ready for handoff" have now been handed back in sync with
# en-US (no changes came back), so:
#   - the status text changes everywhere it is shown (Overview + per
#     language sheets)
#   - the per-language sheets gain a "Latest Target File" / "Latest
#     Handback File" hyperlink (identical to the source / handoff file,
#     since nothing changed) and a new "Latest Handback DateTime"

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- 1. Update the status text everywhere it appears -----------------
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZhCn.Range("B2").Value = $newStatus
$wsZhCn.Range("B3").Value = $newStatus

$wsDeDe.Range("B2").Value = $newStatus
$wsDeDe.Range("B3").Value = $newStatus

# --- 2. zh-cn sheet: fill in Latest Target File / Latest Handback File
#        hyperlinks and Latest Handback DateTime ----------------------
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("E2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/264a6fee652ba6d07ecef8f8256083ebb2b2d7cc/e2e/0684480c-6f3c-425f-a471-6d4b5ac73605.md",
    "",
    "",
    "0684480c-6f3c-425f-a471-6d4b5ac73605.md"
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/55d243506468055bbb36f72f9b3b899a80036449/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/0684480c-6f3c-425f-a471-6d4b5ac73605.dfe660d99174164cc12bfea3eb4ec963f6c27105.zh-cn.xlf",
    "",
    "",
    "0684480c-6f3c-425f-a471-6d4b5ac73605.dfe660d99174164cc12bfea3eb4ec963f6c27105.zh-cn.xlf"
)
$wsZhCn.Range("G2").Value = "2016-02-17 06:11:33"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("E3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/264a6fee652ba6d07ecef8f8256083ebb2b2d7cc/e2e/7bc79fcf-000e-47b3-a00c-ff4582d3354f.md",
    "",
    "",
    "7bc79fcf-000e-47b3-a00c-ff4582d3354f.md"
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/55d243506468055bbb36f72f9b3b899a80036449/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/7bc79fcf-000e-47b3-a00c-ff4582d3354f.eeb370b2f6e8c1b7d8df8d6a9898f6ee4e9852c2.zh-cn.xlf",
    "",
    "",
    "7bc79fcf-000e-47b3-a00c-ff4582d3354f.eeb370b2f6e8c1b7d8df8d6a9898f6ee4e9852c2.zh-cn.xlf"
)
$wsZhCn.Range("G3").Value = "2016-02-17 06:11:50"

# --- 3. de-de sheet: same treatment -----------------------------------
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("E2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/264a6fee652ba6d07ecef8f8256083ebb2b2d7cc/e2e/0684480c-6f3c-425f-a471-6d4b5ac73605.md",
    "",
    "",
    "0684480c-6f3c-425f-a471-6d4b5ac73605.md"
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e4ef128c3eeeebab340f8ed90cc1d3a30433f9a2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/0684480c-6f3c-425f-a471-6d4b5ac73605.dfe660d99174164cc12bfea3eb4ec963f6c27105.de-de.xlf",
    "",
    "",
    "0684480c-6f3c-425f-a471-6d4b5ac73605.dfe660d99174164cc12bfea3eb4ec963f6c27105.de-de.xlf"
)
$wsDeDe.Range("G2").Value = "2016-02-17 06:11:50"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("E3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/264a6fee652ba6d07ecef8f8256083ebb2b2d7cc/e2e/7bc79fcf-000e-47b3-a00c-ff4582d3354f.md",
    "",
    "",
    "7bc79fcf-000e-47b3-a00c-ff4582d3354f.md"
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e4ef128c3eeeebab340f8ed90cc1d3a30433f9a2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/7bc79fcf-000e-47b3-a00c-ff4582d3354f.eeb370b2f6e8c1b7d8df8d6a9898f6ee4e9852c2.de-de.xlf",
    "",
    "",
    "7bc79fcf-000e-47b3-a00c-ff4582d3354f.eeb370b2f6e8c1b7d8df8d6a9898f6ee4e9852c2.de-de.xlf"
)
$wsDeDe.Range("G3").Value = "2016-02-17 06:11:50"
